$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gate Drivers")
$dc = $wb.Worksheets.Item("DC-DC Converter")
$dc.Range("C7").Copy()   # style 8: border1 only, no alignment
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").NumberFormat = "#,##0"

$dc.Range("C7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Style = "Yüzde"     # xfId4 style? need to check name in cellStyles
Write-Host "trying style apply"
